$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$e2 = @'

                    <filter>
                        <interfaces xmlns="http://openconfig.net/yang/interfaces">
                        <interface>
                        <name>1/1/1</name>
                        </interface>
                        </interfaces>
                    </filter>
                    
'@

$f2 = @'
<?xml version="1.0" encoding="UTF-8"?>
<rpc-reply message-id="urn:uuid:02acc40c-f5ec-4807-a70b-93e7288e84a9" xmlns:nc="urn:ietf:params:xml:ns:netconf:base:1.0" xmlns="urn:ietf:params:xml:ns:netconf:base:1.0">
    <data>
        <interfaces xmlns="http://openconfig.net/yang/interfaces">
            <interface>
                <name>1/1/1</name>
                <config>
                    <name>1/1/1</name>
                    <type>ethernetCsmacd</type>
                    <mtu>1500</mtu>
                    <description>test</description>
                    <enabled>true</enabled>
                </config>
                <ethernet xmlns="http://openconfig.net/yang/interfaces/ethernet">
                    <config>
                        <port-speed>SPEED_100MB</port-speed>
                    </config>
                </ethernet>
            </interface>
        </interfaces>
    </data>
</rpc-reply>
'@

$g2 = @'
  <edit-config>
    <target>
      <candidate/>
    </target>
    <config>
      <interfaces xmlns="http://openconfig.net/yang/interfaces">
        <interface>
          <name>1/1/1</name>
          <ethernet xmlns="http://openconfig.net/yang/interfaces/ethernet">
            <config>
              <auto-negotiate>false</auto-negotiate>
            </config>
          </ethernet>
        </interface>
      </interfaces>
    </config>
  </edit-config>
'@

$h2 = @'
- Response of edit-config: <?xml version="1.0" encoding="UTF-8"?>
<rpc-reply message-id="urn:uuid:360eb500-4805-4b58-b796-24ae9ec4f348" xmlns:nc="urn:ietf:params:xml:ns:netconf:base:1.0" xmlns="urn:ietf:params:xml:ns:netconf:base:1.0">
    <ok/>
</rpc-reply> 
 - Response of commit: <?xml version="1.0" encoding="UTF-8"?>
<rpc-reply message-id="urn:uuid:284458d7-521c-4d95-a832-500da58b72f7" xmlns:nc="urn:ietf:params:xml:ns:netconf:base:1.0" xmlns="urn:ietf:params:xml:ns:netconf:base:1.0">
    <ok/>
</rpc-reply>
'@

$i2 = @'
<?xml version="1.0" encoding="UTF-8"?>
<rpc-reply message-id="urn:uuid:52237c29-63fa-4f12-b2bf-6759fc20fb24" xmlns:nc="urn:ietf:params:xml:ns:netconf:base:1.0" xmlns="urn:ietf:params:xml:ns:netconf:base:1.0">
    <data>
        <interfaces xmlns="http://openconfig.net/yang/interfaces">
            <interface>
                <name>1/1/1</name>
                <config>
                    <name>1/1/1</name>
                    <type>ethernetCsmacd</type>
                    <mtu>1500</mtu>
                    <description>test</description>
                    <enabled>true</enabled>
                </config>
                <ethernet xmlns="http://openconfig.net/yang/interfaces/ethernet">
                    <config>
                        <auto-negotiate>false</auto-negotiate>
                        <port-speed>SPEED_100MB</port-speed>
                    </config>
                </ethernet>
            </interface>
        </interfaces>
    </data>
</rpc-reply>
'@

$ws.Range("E2").Value = $e2
$ws.Range("F2").Value = $f2
$ws.Range("G2").Value = $g2
$ws.Range("H2").Value = $h2
$ws.Range("I2").Value = $i2
